# Applies the "Anonimyzed fedcore" update:
#  - rename the "fedcore" column header(s) to "approach" on both sheets
#  - add a top+bottom border under the B1:D1 / E1:G1 merged header cells,
#    with an extra right-hand border on the last column of each group
#  - drop the stray empty inline-string cell at G5 on the second sheet

$wb  = $excel.ActiveWorkbook
$app = $excel

# Border/paste constants (xlEdgeTop/Bottom/Right, xlContinuous, xlPasteFormats)
$xlEdgeTop      = 8
$xlEdgeBottom   = 9
$xlEdgeRight    = 10
$xlContinuous   = 1
$xlPasteFormats = -4122

$ws1 = $wb.Worksheets.Item(1)   # quality_comparison
$ws2 = $wb.Worksheets.Item(2)   # computational_comparison

# Build the two new border styles exactly once (on sheet 1), then replicate
# them everywhere else via copy / paste-special so the workbook's shared
# style table only ever gains the two new entries it needs.
$c1 = $ws1.Range("C1")
$c1.ClearFormats()
$c1.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
$c1.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous

$d1 = $ws1.Range("D1")
$d1.ClearFormats()
$d1.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
$d1.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
$d1.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous

# Propagate the "top+bottom" style to the other cells that need it
$c1.Copy()
$ws2.Range("C1").PasteSpecial($xlPasteFormats)
$ws2.Range("F1").PasteSpecial($xlPasteFormats)

# Propagate the "top+bottom+right" style to the other cells that need it
$d1.Copy()
$ws2.Range("D1").PasteSpecial($xlPasteFormats)
$ws2.Range("G1").PasteSpecial($xlPasteFormats)

$app.CutCopyMode = 0

# ---- Text updates: anonymize "fedcore" -> "approach" ----
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# ---- Drop the stray empty inline-string cell ----
$ws2.Range("G5").ClearContents()
